$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (row 1 is the header), shifting the existing
# data rows (2-23) down to (3-24).
$ws.Range("A2").EntireRow.Insert()

# The inserted row inherits formatting from the row above (the bold header);
# strip that back to plain formatting like every other data row, then give
# the date cell (D2) the same date number format used by the rest of
# column D.
$ws.Range("A2:T2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with this week's entry (same market/product
# template as the rest of the sheet, new date + price figures).
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Femacal de La Calera"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44616
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104001
$ws.Range("J2").Value = "Granada"
$ws.Range("K2").Value = "Wonderfull"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("Q2").Value = "$/caja 14 kilos empedrada"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 14

